# Keyboard shortcuts workbook update
# - CRTL+N (A7) becomes a tri-run rich string: "COMMAND-N" + nbsp-spacer + "[Windows: CRTL+N]"
# - ALT+h... (A21) becomes OPTION+h... (simple text swap, keeps its existing formatting)
# - CRTL 0-9 (A26) becomes a tri-run rich string: "COMMAND 0-9" + nbsp-spacer + "[Windows: CRTL 0-9]"
# - Rows 7 and 26 shrink from 15pt to 14.6pt
# - Selection moves to A26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blue = 0xC07000   # BGR for RGB 00,70,C0 (FF0070C0 incl. alpha)
$nbsp3 = '\u00A0\u00A0\u00A0'

# ---- A7 : CRTL+N -> COMMAND-N <nbsp> [Windows: CRTL+N] ----
$cellA7 = $ws.Range("A7")
$run1 = "COMMAND-N"
$run2 = $nbsp3
$run3 = "[Windows: CRTL+N]"
$cellA7.Value = $run1 + $run2 + $run3
$cellA7.Characters(1, $run1.Length).Font.Color = $blue
$cellA7.Characters($run1.Length + 1, $run2.Length).Font.Color = $blue
$cellA7.Characters($run1.Length + $run2.Length + 1, $run3.Length).Font.Color = $blue

# ---- A21 : ALT+h... -> OPTION+h... ----
$ws.Range("A21").Value = "OPTION+h" + $nbsp3 + "[Windows: CTRL+SHIFT+h]"

# ---- A26 : CRTL 0-9 -> COMMAND 0-9 <nbsp> [Windows: CRTL 0-9] ----
$cellA26 = $ws.Range("A26")
$run1b = "COMMAND 0-9"
$run2b = $nbsp3
$run3b = "[Windows: CRTL 0-9]"
$cellA26.Value = $run1b + $run2b + $run3b
$cellA26.Characters(1, $run1b.Length).Font.Color = $blue
$cellA26.Characters($run1b.Length + 1, $run2b.Length).Font.Color = $blue
$cellA26.Characters($run1b.Length + $run2b.Length + 1, $run3b.Length).Font.Color = $blue

# ---- Row heights ----
$ws.Rows.Item(7).RowHeight = 14.6
$ws.Rows.Item(26).RowHeight = 14.6

# ---- Default column width (cosmetic) ----
$ws.StandardWidth = 8.578125

# ---- Selection ----
$ws.Range("A26").Select() | Out-Null
